# ProductBacklog.xlsx update
# - Fix typo in SP9 description
# - Rename "Done" column (F) to "Logic" and add three new tracking
#   columns: Design (G), Tested (H), Done (I)
# - Fill the new/changed boolean columns per backlog item

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "hist" -> "his"
$ws.Range("C8").Value = "User wants to review his tutor"

# Header row: rename column F and add new headers G:I, copying the
# bold/shaded header format from F1 onto the new header cells.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:I1").PasteSpecial(-4122) | Out-Null

$ws.Range("F1").Value = "Logic"
$ws.Range("G1").Value = "Design"
$ws.Range("H1").Value = "Tested"
$ws.Range("I1").Value = "Done"

# Data rows 2-14: new values for F (Logic) and the newly added
# G (Design), H (Tested), I (Done) columns.
$values = @{
    2  = @("no",  "no",  "no",  "no")
    3  = @("yes", "no",  "no",  "no")
    4  = @("no",  "no",  "no",  "no")
    5  = @("yes", "no",  "no",  "no")
    6  = @("no",  "no",  "no",  "no")
    7  = @("yes", "no",  "no",  "no")
    8  = @("no",  "no",  "no",  "no")
    9  = @("no",  "no",  "no",  "no")
    10 = @("yes", "no",  "no",  "no")
    11 = @("no",  "no",  "no",  "no")
    12 = @("yes", "yes", "no",  "yes")
    13 = @("yes", "no",  "no",  "no")
    14 = @("yes", "yes", "no",  "yes")
}

foreach ($r in 2..14) {
    $row = $values[$r]
    $ws.Cells.Item($r, 6).Value = $row[0]
    $ws.Cells.Item($r, 7).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
    $ws.Cells.Item($r, 9).Value = $row[3]
}
